# This script applies a permutation of the data rows (rows 2-11, columns A:AY)
# in the active worksheet, as described by the target diff. The row numbers
# themselves do not move, but the record data contained in them gets
# reshuffled into a different order (a pure data permutation, e.g. two
# duplicate/near-duplicate observation rows effectively being renumbered).
#
# Mapping (new row -> row the data currently comes from):
#   2  <- 5
#   3  <- 2
#   4  <- 3
#   5  <- 4
#   6  <- 10
#   7  <- 9
#   8  <- 7
#   9  <- 6
#   10 <- 11
#   11 <- 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol  = "AY"

# Snapshot all the source rows first, so that writing results doesn't
# clobber data that is still needed for a later step of the permutation
# (the mapping contains cycles, e.g. 6 -> 10 -> 11 -> 8 -> 7 -> 9 -> 6).
$rowData = @{}
foreach ($r in 2..11) {
    $rng = $ws.Range("$firstCol$r`:$lastCol$r")
    $rowData[$r] = $rng.Value2
}

$mapping = @{
    2  = 5
    3  = 2
    4  = 3
    5  = 4
    6  = 10
    7  = 9
    8  = 7
    9  = 6
    10 = 11
    11 = 8
}

# Columns Y, Z, AA, AB hold date/time values stored as plain text
# (e.g. "2023-08-29", "00:00"). Writing such strings back through
# Value2 would otherwise be auto-parsed by Excel into date/time serial
# numbers, so force those destination columns to Text format first so
# the values round-trip as the original strings.
foreach ($newRow in 2..11) {
    $ws.Range("Y$newRow`:AB$newRow").NumberFormat = "@"
}

foreach ($newRow in 2..11) {
    $srcRow = $mapping[$newRow]
    $destRng = $ws.Range("$firstCol$newRow`:$lastCol$newRow")
    $destRng.Value2 = $rowData[$srcRow]
}
